$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- Update panel_query_time values (column F) in the "data" sheet ---
$ws.Cells.Item(2, 6).Value = "2021-10-05 14:34:18.246243"
$ws.Cells.Item(3, 6).Value = "2021-10-05 14:34:18.246251"
$ws.Cells.Item(4, 6).Value = "2021-10-05 14:34:18.246254"
$ws.Cells.Item(5, 6).Value = "2021-10-05 14:34:18.246257"
$ws.Cells.Item(6, 6).Value = "2021-10-05 14:34:18.246260"
$ws.Cells.Item(7, 6).Value = "2021-10-05 14:34:18.246262"
$ws.Cells.Item(8, 6).Value = "2021-10-05 14:34:18.246265"
$ws.Cells.Item(9, 6).Value = "2021-10-05 14:34:18.246267"
$ws.Cells.Item(10, 6).Value = "2021-10-05 14:34:18.246270"
$ws.Cells.Item(11, 6).Value = "2021-10-05 14:34:18.246272"
$ws.Cells.Item(12, 6).Value = "2021-10-05 14:34:18.246275"
$ws.Cells.Item(13, 6).Value = "2021-10-05 14:34:18.246277"
$ws.Cells.Item(14, 6).Value = "2021-10-05 14:34:18.246280"
$ws.Cells.Item(15, 6).Value = "2021-10-05 14:34:18.246282"
$ws.Cells.Item(16, 6).Value = "2021-10-05 14:34:18.246285"
$ws.Cells.Item(17, 6).Value = "2021-10-05 14:34:18.246287"
$ws.Cells.Item(18, 6).Value = "2021-10-05 14:34:18.246290"
$ws.Cells.Item(19, 6).Value = "2021-10-05 14:34:18.246292"
$ws.Cells.Item(20, 6).Value = "2021-10-05 14:34:18.246295"
$ws.Cells.Item(21, 6).Value = "2021-10-05 14:34:18.246297"
$ws.Cells.Item(22, 6).Value = "2021-10-05 14:34:18.246300"
$ws.Cells.Item(23, 6).Value = "2021-10-05 14:34:18.246302"
$ws.Cells.Item(24, 6).Value = "2021-10-05 14:34:18.246305"
$ws.Cells.Item(25, 6).Value = "2021-10-05 14:34:18.246307"
$ws.Cells.Item(26, 6).Value = "2021-10-05 14:34:18.246310"
$ws.Cells.Item(27, 6).Value = "2021-10-05 14:34:18.246313"
$ws.Cells.Item(28, 6).Value = "2021-10-05 14:34:18.246315"
$ws.Cells.Item(29, 6).Value = "2021-10-05 14:34:18.246318"
$ws.Cells.Item(30, 6).Value = "2021-10-05 14:34:18.246320"
$ws.Cells.Item(31, 6).Value = "2021-10-05 14:34:18.246323"
$ws.Cells.Item(32, 6).Value = "2021-10-05 14:34:18.246325"
$ws.Cells.Item(33, 6).Value = "2021-10-05 14:34:18.246328"
$ws.Cells.Item(34, 6).Value = "2021-10-05 14:34:18.246331"
$ws.Cells.Item(35, 6).Value = "2021-10-05 14:34:18.246333"
$ws.Cells.Item(36, 6).Value = "2021-10-05 14:34:18.246336"
$ws.Cells.Item(37, 6).Value = "2021-10-05 14:34:18.246338"
$ws.Cells.Item(38, 6).Value = "2021-10-05 14:34:18.246341"
$ws.Cells.Item(39, 6).Value = "2021-10-05 14:34:18.246343"
$ws.Cells.Item(40, 6).Value = "2021-10-05 14:34:18.246346"
$ws.Cells.Item(41, 6).Value = "2021-10-05 14:34:18.246348"
$ws.Cells.Item(42, 6).Value = "2021-10-05 14:34:18.246351"
$ws.Cells.Item(43, 6).Value = "2021-10-05 14:34:18.246354"
$ws.Cells.Item(44, 6).Value = "2021-10-05 14:34:18.246356"
$ws.Cells.Item(45, 6).Value = "2021-10-05 14:34:18.246359"
$ws.Cells.Item(46, 6).Value = "2021-10-05 14:34:18.246361"
$ws.Cells.Item(47, 6).Value = "2021-10-05 14:34:18.246364"
$ws.Cells.Item(48, 6).Value = "2021-10-05 14:34:18.246366"
$ws.Cells.Item(49, 6).Value = "2021-10-05 14:34:18.246369"
$ws.Cells.Item(50, 6).Value = "2021-10-05 14:34:18.246371"
$ws.Cells.Item(51, 6).Value = "2021-10-05 14:34:18.246373"
$ws.Cells.Item(52, 6).Value = "2021-10-05 14:34:18.246376"
$ws.Cells.Item(53, 6).Value = "2021-10-05 14:34:18.246378"
$ws.Cells.Item(54, 6).Value = "2021-10-05 14:34:18.246381"
$ws.Cells.Item(55, 6).Value = "2021-10-05 14:34:18.246384"
$ws.Cells.Item(56, 6).Value = "2021-10-05 14:34:18.246387"
$ws.Cells.Item(57, 6).Value = "2021-10-05 14:34:18.246389"
$ws.Cells.Item(58, 6).Value = "2021-10-05 14:34:18.246391"
$ws.Cells.Item(59, 6).Value = "2021-10-05 14:34:18.246394"
$ws.Cells.Item(60, 6).Value = "2021-10-05 14:34:18.246396"
$ws.Cells.Item(61, 6).Value = "2021-10-05 14:34:18.246399"
$ws.Cells.Item(62, 6).Value = "2021-10-05 14:34:18.246402"
$ws.Cells.Item(63, 6).Value = "2021-10-05 14:34:18.246404"
$ws.Cells.Item(64, 6).Value = "2021-10-05 14:34:18.246407"
$ws.Cells.Item(65, 6).Value = "2021-10-05 14:34:18.246409"
$ws.Cells.Item(66, 6).Value = "2021-10-05 14:34:18.246413"
$ws.Cells.Item(67, 6).Value = "2021-10-05 14:34:18.246416"
$ws.Cells.Item(68, 6).Value = "2021-10-05 14:34:18.246418"
$ws.Cells.Item(69, 6).Value = "2021-10-05 14:34:18.246421"
$ws.Cells.Item(70, 6).Value = "2021-10-05 14:34:18.246423"
$ws.Cells.Item(71, 6).Value = "2021-10-05 14:34:18.246426"
$ws.Cells.Item(72, 6).Value = "2021-10-05 14:34:18.246429"
$ws.Cells.Item(73, 6).Value = "2021-10-05 14:34:18.246431"
$ws.Cells.Item(74, 6).Value = "2021-10-05 14:34:18.246434"
$ws.Cells.Item(75, 6).Value = "2021-10-05 14:34:18.246436"
$ws.Cells.Item(76, 6).Value = "2021-10-05 14:34:18.246439"
$ws.Cells.Item(77, 6).Value = "2021-10-05 14:34:18.246441"
$ws.Cells.Item(78, 6).Value = "2021-10-05 14:34:18.246446"
$ws.Cells.Item(79, 6).Value = "2021-10-05 14:34:18.246449"
$ws.Cells.Item(80, 6).Value = "2021-10-05 14:34:18.246452"
$ws.Cells.Item(81, 6).Value = "2021-10-05 14:34:18.246454"
$ws.Cells.Item(82, 6).Value = "2021-10-05 14:34:18.246457"
$ws.Cells.Item(83, 6).Value = "2021-10-05 14:34:18.246460"
$ws.Cells.Item(84, 6).Value = "2021-10-05 14:34:18.246462"
$ws.Cells.Item(85, 6).Value = "2021-10-05 14:34:18.246465"
$ws.Cells.Item(86, 6).Value = "2021-10-05 14:34:18.246467"
$ws.Cells.Item(87, 6).Value = "2021-10-05 14:34:18.246470"
$ws.Cells.Item(88, 6).Value = "2021-10-05 14:34:18.246472"

# --- Add the new "metadata" sheet after "data" ---
$newSheet = $wb.Worksheets.Add($null, $ws)
$newSheet.Name = "metadata"

# Header row (bold, bordered, centered - matches the "data" sheet's header style)
$newSheet.Range("B1").Value = "data_name"
$newSheet.Range("C1").Value = "data_id"
$newSheet.Range("D1").Value = "data_version"
$newSheet.Range("E1").Value = "data_version_created"
$newSheet.Range("F1").Value = "panel_query_time"
$newSheet.Range("G1").Value = "panel_get_request"

$headerRng = $newSheet.Range("B1:G1")
$headerRng.Font.Bold = $true
$headerRng.Borders.LineStyle = 1
$headerRng.HorizontalAlignment = -4108
$headerRng.VerticalAlignment = -4160

# Data row 2
$newSheet.Range("A2").Value = 0
$idxRng = $newSheet.Range("A2")
$idxRng.Font.Bold = $true
$idxRng.Borders.LineStyle = 1
$idxRng.HorizontalAlignment = -4108
$idxRng.VerticalAlignment = -4160

$newSheet.Range("B2").Value = "Leukodystrophy - adult onset"
$newSheet.Range("C2").Value = 299
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "0.92"
$newSheet.Range("D2").Style = "Normal"
$newSheet.Range("E2").Value = "2021-09-16T08:20:42.902711Z"
$newSheet.Range("F2").Value = "2021-10-05 14:34:18.242955"
$newSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/299/?format=json"

Write-Output "metadata sheet added"
